# Add a new "Connector" worksheet (USB parts) after "Transistor",
# populate it with the USB connector comparison table, and apply the
# same visual formatting conventions used on the "Transistor" sheet.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Transistor")

# --- create + position the new sheet -------------------------------------
$ws = $wb.Worksheets.Add($null, $src)
$ws.Name = "Connector"

# --- column widths (approximate the source workbook's hand-tuned widths) --
$ws.Columns.Item(1).ColumnWidth = 60.17
$ws.Columns.Item(2).ColumnWidth = 38.02
$ws.Columns.Item(3).ColumnWidth = 38.04
$ws.Columns.Item(4).ColumnWidth = 38.04
$ws.Columns.Item(5).ColumnWidth = 38.04
$ws.Columns.Item(6).ColumnWidth = 38.04
$ws.Columns.Item(7).ColumnWidth = 38.04
$ws.Columns.Item(8).ColumnWidth = 38.04

# --- row heights ------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 23
for ($r = 2; $r -le 10; $r++) {
    $ws.Rows.Item($r).RowHeight = 26
}

# --- cell values (8 columns x 10 rows) --------------------------------------
$rows = @(
    @('USB','47346-0001','10033526-N3212MLF','48037-2200','UE27-AC54-100','UE27-AE54-100','87583-2010BLF','61729-0010BLF'),
    @('Manufacturer','MOLEX','Amphenol FCI','MOLEX','Amphenol','Amphenol','Amphenol FCI','Amphenol FCI'),
    @('USB Type','Micro USB Type B Receptacle','Mini USB Type B Receptacle','USB Type A Plug','USB Type A Receptacle','USB Type A Receptacle','USB Type A Receptacle','USB Type B Receptacle'),
    @('Gender','Female','Female','Male','Female','Female','Female','Female'),
    @('Standard','USB 2.0','USB 2.0','USB','USB 2.0','USB 2.0','USB 2.0','USB 2.0'),
    @('Current Rating ','1.8 A ','-','1.5A','-','-','-','-'),
    @('Voltage Rating ','30 VAC','-','150 VAC','-','-','-','-'),
    @(('Operating Temperature'+[char]0x00A0),'- 20 C to + 85 C','-','- 20 C to + 85 C','-','-','-','-'),
    @('Mounting Type','Surface Mount','Surface Mount','Surface Mount','Through Hole','Through Hole','Surface Mount','Through Hole'),
    @('Mounting Angle','Right','Right','Right','Right','Vertical','Right','Right')
)

# Style catalogue, keyed the same way the source sheet uses them:
#   2 -> header band (bold white on red)         e.g. Transistor!A1
#   3 -> label cell, blue fill, bold              e.g. Transistor!A2
#   4 -> value cell, blue fill                    e.g. Transistor!B2
#   5 -> label cell, no fill, bold                e.g. Transistor!A3
#   6 -> value cell, no fill                      e.g. Transistor!B3
#   9 -> value cell, blue fill, wrapped + left     (4, plus WrapText)
$styleIdx = @(
    @(2,2,2,2,2,2,2,2),
    @(3,4,9,4,9,9,9,9),
    @(5,6,6,6,6,6,6,6),
    @(3,4,4,4,4,4,4,4),
    @(5,6,6,6,6,6,6,6),
    @(3,4,4,4,4,4,4,4),
    @(5,6,6,6,6,6,6,6),
    @(3,4,4,4,4,4,4,4),
    @(5,6,6,6,6,6,6,6),
    @(3,4,4,4,4,4,4,4)
)

$fmtSrc = @{
    2 = $src.Range("A1")
    3 = $src.Range("A2")
    4 = $src.Range("B2")
    5 = $src.Range("A3")
    6 = $src.Range("B3")
    9 = $src.Range("B2")
}

# Copy each distinct source format once, then paint it onto every
# destination cell that needs it (cheap: avoids re-copying per cell).
foreach ($key in $fmtSrc.Keys) {
    $fmtSrc[$key].Copy()
    for ($r = 1; $r -le 10; $r++) {
        for ($c = 1; $c -le 8; $c++) {
            if ($styleIdx[$r-1][$c-1] -eq $key) {
                $ws.Cells.Item($r, $c).PasteSpecial(-4122) | Out-Null
            }
        }
    }
}

# Write values + finish the wrap/left-align tweak for style 9 cells.
for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $rows[$r-1][$c-1]
        if ($styleIdx[$r-1][$c-1] -eq 9) {
            $cell.HorizontalAlignment = -4131   # xlLeft
            $cell.WrapText = $true
        }
    }
}

$excel.CutCopyMode = $false

# --- sheet-level view / page setup, mirroring the Transistor sheet ---------
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false

$ws.PageSetup.LeftMargin   = 54
$ws.PageSetup.RightMargin  = 54
$ws.PageSetup.TopMargin    = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Zoom = 100
$ws.PageSetup.Orientation = 1
$ws.PageSetup.CenterFooter = '&"Helvetica,Regular"&12&K000000&P'

# restore the original active sheet
$src.Activate()

Write-Host "Connector sheet added with $($ws.UsedRange.Rows.Count) rows x $($ws.UsedRange.Columns.Count) cols"
